$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.921.13"
$ws.Range("D3").Value = "1.875.36"
$ws.Range("E3").Value = "  -0.94%  "
$c = $ws.Range("D4")
$c.NumberFormat = "@"
$c.Value = "0.9994"
$c.Style = "Normal"
$ws.Range("E4").Value = "  -0.12%  "
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "0.7412"
$c.Style = "Normal"
$ws.Range("E5").Value = "  -3.98%  "
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "242.49"
$c.Style = "Normal"
$ws.Range("E6").Value = "  -0.65%  "
$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = "1.000"
$c.Style = "Normal"
$ws.Range("E7").Value = "  -0.06%  "
$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = "0.3154"
$c.Style = "Normal"
$ws.Range("E8").Value = "  +0.89%  "
$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = "0.07177"
$c.Style = "Normal"
$ws.Range("E9").Value = "  -0.72%  "
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = "24.76"
$c.Style = "Normal"
$ws.Range("E10").Value = "  -3.53%  "
$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = "0.08449"
$c.Style = "Normal"
$ws.Range("E11").Value = "  -2.88%  "
$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = "0.7536"
$c.Style = "Normal"
$ws.Range("E12").Value = "  -2.29%  "
$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").Value = "1.893.74"
$ws.Range("E13").Value = "  -7.95%  "
$ws.Range("B14").Value = "Polkadot"
$ws.Range("C14").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = "5.401"
$c.Style = "Normal"
$ws.Range("E14").Value = "  -0.08%  "
$ws.Range("E15").Value = "  -1.62%  "
$ws.Range("D16").Value = "29.930.02"
$ws.Range("E16").Value = "  -0.39%  "
$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = "6.103"
$c.Style = "Normal"
$ws.Range("E17").Value = "  -1.72%  "
$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = "13.60"
$c.Style = "Normal"
$ws.Range("E18").Value = "  -2.33%  "
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = "243.58"
$c.Style = "Normal"
$ws.Range("E19").Value = "  -0.66%  "
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = "0.000007824"
$c.Style = "Normal"
$ws.Range("E20").Value = "  -0.40%  "
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = "0.9998"
$c.Style = "Normal"
$ws.Range("E21").Value = "  -0.23%  "
$ws.Range("D22").Value = "2.124.70"
$ws.Range("E22").Value = "  -8.98%  "
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = "8.008"
$c.Style = "Normal"
$ws.Range("E23").Value = "  -2.39%  "
$ws.Range("E24").Value = "  -0.23%  "
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = "0.1561"
$c.Style = "Normal"
$ws.Range("E25").Value = "  -1.99%  "
$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = "9.336"
$c.Style = "Normal"
$ws.Range("E26").Value = "  -1.85%  "
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = "165.82"
$c.Style = "Normal"
$ws.Range("E27").Value = "  +2.14%  "
$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = "18.64"
$c.Style = "Normal"
$ws.Range("E28").Value = "  -0.88%  "
$ws.Range("E29").Value = "  +0.11%  "
$ws.Range("E30").Value = "  +2.95%  "
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = "4.617"
$c.Style = "Normal"
$ws.Range("E31").Value = "  +2.26%  "
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = "1.530"
$c.Style = "Normal"
$ws.Range("E32").Value = "  -0.95%  "
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = "4.282"
$c.Style = "Normal"
$ws.Range("E33").Value = "  +4.00%  "
$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = "0.05339"
$c.Style = "Normal"
$ws.Range("E34").Value = "  -2.36%  "
$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = "1.244"
$c.Style = "Normal"
$ws.Range("E35").Value = "  -0.38%  "
$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = "0.7566"
$c.Style = "Normal"
$ws.Range("E36").Value = "  +0.47%  "
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = "0.9999"
$c.Style = "Normal"
$ws.Range("E37").Value = "  -0.42%  "
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = "2.698"
$c.Style = "Normal"
$ws.Range("E38").Value = "  -0.23%  "
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = "0.01959"
$c.Style = "Normal"
$ws.Range("E39").Value = "  -0.58%  "
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = "2.750"
$c.Style = "Normal"
$ws.Range("E40").Value = "  -1.33%  "
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = "0.4494"
$c.Style = "Normal"
$ws.Range("E41").Value = "  -0.33%  "
$ws.Range("D42").Value = "1.113.58"
$ws.Range("E42").Value = "  +1.57%  "
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = "6.097"
$c.Style = "Normal"
$ws.Range("E43").Value = "  +0.80%  "
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = "72.50"
$c.Style = "Normal"
$ws.Range("E44").Value = "  -1.49%  "
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = "0.8600"
$c.Style = "Normal"
$ws.Range("E45").Value = "  +0.69%  "
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = "103.27"
$c.Style = "Normal"
$ws.Range("E47").Value = "  +0.42%  "
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = "7.686"
$c.Style = "Normal"
$ws.Range("E48").Value = "  +0.81%  "
$ws.Range("B49").Value = "RenderToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = "1.844"
$c.Style = "Normal"
$ws.Range("E49").Value = "  -2.03%  "
$ws.Range("B50").Value = "SynthetixNetwork"
$ws.Range("C50").Value = "https://coinranking.com/coin/sgxZRXbK0FDc+synthetixnetwork-snx"
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = "3.069"
$c.Style = "Normal"
$ws.Range("E50").Value = "  +3.03%  "
$ws.Range("D51").Value = "2.022.66"
$ws.Range("E51").Value = "  -8.32%  "
